# Documentation - Sewage - Updated descriptive text and installed and demand
# load calculations to new floor plant.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated load calculation figures (col D = unit power, col E = quantity,
#     col F = total power, computed via shared formula E*D) ---

# Lampada Compacta: quantity 15 -> 9 (row 3)
$ws.Range("E3").Value = 9

# Misturador Submerso: unit power 2.65 -> 2.66 (row 6)
$ws.Range("D6").Value = 2.66

# Bomba Helicoidal: unit power 3.5 -> 3.55 (row 8)
$ws.Range("D8").Value = 3.55

# Bomba Centrifuga Submersivel: unit power 10.4 -> 10 (row 9)
$ws.Range("D9").Value = 10

# Soprador: unit power 0.93 -> 0.31 (row 10)
$ws.Range("D10").Value = 0.31

# --- Move the saved selection cursor to reflect where the author left off ---
$ws.Range("B18").Select() | Out-Null
